$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 82.73
$ws.Range("I15").Value = 82.73
$ws.Range("K15").Value = 248.19
$ws.Range("M15").Value = -79.19
$ws.Range("H59").Value = 1125
$ws.Range("J59").Value = 1125
$ws.Range("L59").Value = 3375
$ws.Range("N59").Value = -4489
$ws.Range("H124").Value = 43780
$ws.Range("J124").Value = 43780
$ws.Range("L124").Value = 43780
$ws.Range("N124").Value = -53600
$ws.Range("H129").Value = 909.45
$ws.Range("I129").Value = 359.4
$ws.Range("J129").Value = 959.4545000000001
$ws.Range("K129").Value = 1078.2
$ws.Range("L129").Value = 2878.3635
$ws.Range("M129").Value = 3921.8
$ws.Range("N129").Value = -12878.3635
$ws.Range("H137").Value = 2887.394
$ws.Range("I137").Value = 2059.75
$ws.Range("K137").Value = 6179.25
$ws.Range("M137").Value = -3629.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 39600
$ws.Range("J7").Value = 39600
$ws.Range("L7").Value = 39600
$ws.Range("N7").Value = -39828
$ws.Range("H42").Value = 32499
$ws.Range("J42").Value = 32499
$ws.Range("L42").Value = 32499
$ws.Range("N42").Value = -33471
$ws.Range("H61").Value = 2434.5454
$ws.Range("I61").Value = 1140
$ws.Range("J61").Value = 3513.3333
$ws.Range("K61").Value = 1140
$ws.Range("L61").Value = 3513.3333
$ws.Range("M61").Value = -928
$ws.Range("N61").Value = -3937.3333
$ws.Range("H74").Value = 1369.1482
$ws.Range("I74").Value = 978.2292
$ws.Range("J74").Value = 4496.5
$ws.Range("K74").Value = 978.2292
$ws.Range("L74").Value = 4496.5
$ws.Range("M74").Value = -104.2292
$ws.Range("N74").Value = -6244.5
$ws.Range("H77").Value = 1369.1482
$ws.Range("I77").Value = 978.2292
$ws.Range("J77").Value = 4496.5
$ws.Range("K77").Value = 4891.146
$ws.Range("L77").Value = 22482.5
$ws.Range("M77").Value = -523.1459999999997
$ws.Range("N77").Value = -31218.5
$ws.Range("H97").Value = 1426
$ws.Range("I97").Value = 1168.6
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 1168.6
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -672.5999999999999
$ws.Range("N97").Value = -4992
$ws.Range("H136").Value = 2434.5454
$ws.Range("I136").Value = 1140
$ws.Range("J136").Value = 3513.3333
$ws.Range("K136").Value = 3420
$ws.Range("L136").Value = 10539.9999
$ws.Range("M136").Value = -870
$ws.Range("N136").Value = -15639.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8055.773
$ws.Range("I20").Value = 1531.2858
$ws.Range("J20").Value = 19473.625
$ws.Range("K20").Value = 1531.2858
$ws.Range("L20").Value = 19473.625
$ws.Range("M20").Value = -1284.2858
$ws.Range("N20").Value = -19967.625
$ws.Range("H86").Value = 1885.4445
$ws.Range("I86").Value = 1769.8572
$ws.Range("J86").Value = 2290
$ws.Range("K86").Value = 1769.8572
$ws.Range("L86").Value = 2290
$ws.Range("M86").Value = -646.8571999999999
$ws.Range("N86").Value = -4536
$ws.Range("H89").Value = 1885.4445
$ws.Range("I89").Value = 1769.8572
$ws.Range("J89").Value = 2290
$ws.Range("K89").Value = 8849.286
$ws.Range("L89").Value = 11450
$ws.Range("M89").Value = -3233.286
$ws.Range("N89").Value = -22682
$ws.Range("H99").Value = 2658.3547
$ws.Range("I99").Value = 1025.2354
$ws.Range("K99").Value = 1025.2354
$ws.Range("M99").Value = 472.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1500
$ws.Range("I15").Value = 1500
$ws.Range("K15").Value = 1500
$ws.Range("M15").Value = -1330
$ws.Range("H31").Value = 3878.1428
$ws.Range("I31").Value = 2057
$ws.Range("J31").Value = 4740.7896
$ws.Range("K31").Value = 2057
$ws.Range("L31").Value = 4740.7896
$ws.Range("M31").Value = -1762
$ws.Range("N31").Value = -5330.7896
$ws.Range("H34").Value = 3878.1428
$ws.Range("I34").Value = 2057
$ws.Range("J34").Value = 4740.7896
$ws.Range("K34").Value = 2057
$ws.Range("L34").Value = 4740.7896
$ws.Range("M34").Value = -1855
$ws.Range("N34").Value = -5144.7896
$ws.Range("H52").Value = 60440
$ws.Range("J52").Value = 60440
$ws.Range("L52").Value = 60440
$ws.Range("N52").Value = -61028
$ws.Range("H58").Value = 2077.2454
$ws.Range("I58").Value = 1828.7307
$ws.Range("K58").Value = 1828.7307
$ws.Range("M58").Value = -1625.7307
$ws.Range("H74").Value = 44650.25
$ws.Range("J74").Value = 44650.25
$ws.Range("L74").Value = 44650.25
$ws.Range("N74").Value = -46398.25
$ws.Range("H77").Value = 44650.25
$ws.Range("J77").Value = 44650.25
$ws.Range("L77").Value = 133950.75
$ws.Range("N77").Value = -142686.75
$ws.Range("H105").Value = 1944.875
$ws.Range("I105").Value = 1374
$ws.Range("K105").Value = 1374
$ws.Range("M105").Value = 373
$ws.Range("H106").Value = 33750
$ws.Range("J106").Value = 33750
$ws.Range("L106").Value = 33750
$ws.Range("N106").Value = -36274
$ws.Range("H122").Value = 2355.7273
$ws.Range("I122").Value = 1725.7059
$ws.Range("K122").Value = 5177.1177
$ws.Range("M122").Value = -2727.1177
$ws.Range("H132").Value = 2889.5789
$ws.Range("I132").Value = 1661.6
$ws.Range("J132").Value = 7494.5
$ws.Range("K132").Value = 4984.799999999999
$ws.Range("L132").Value = 22483.5
$ws.Range("M132").Value = -2454.799999999999
$ws.Range("N132").Value = -27543.5
$ws.Range("H136").Value = 2077.2454
$ws.Range("I136").Value = 1828.7307
$ws.Range("K136").Value = 5486.1921
$ws.Range("M136").Value = -2936.1921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 941.087
$ws.Range("I113").Value = 744.4737
$ws.Range("K113").Value = 2233.4211
$ws.Range("M113").Value = -63.42110000000002
$ws.Range("H132").Value = 3279.2632
$ws.Range("I132").Value = 1374.75
$ws.Range("J132").Value = 3787.1333
$ws.Range("K132").Value = 12372.75
$ws.Range("L132").Value = 34084.1997
$ws.Range("M132").Value = -9842.75
$ws.Range("N132").Value = -39144.1997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1887
$ws.Range("N6").ClearContents()
$ws.Range("H15").Value = 34336.184
$ws.Range("J15").Value = 34336.184
$ws.Range("L15").Value = 34336.184
$ws.Range("N15").Value = -34912.184
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1750
$ws.Range("N16").ClearContents()
$ws.Range("H80").Value = 83335336
$ws.Range("I80").Value = 250000000
$ws.Range("K80").Value = 250000000
$ws.Range("M80").Value = -249999002
$ws.Range("H81").Value = 34336.184
$ws.Range("J81").Value = 34336.184
$ws.Range("L81").Value = 34336.184
$ws.Range("N81").Value = -36332.184
$ws.Range("H83").Value = 83335336
$ws.Range("I83").Value = 250000000
$ws.Range("K83").Value = 1250000000
$ws.Range("M83").Value = -1249995008
$ws.Range("H84").Value = 34336.184
$ws.Range("J84").Value = 34336.184
$ws.Range("L84").Value = 103008.552
$ws.Range("N84").Value = -112992.552
$ws.Range("H97").Value = 2945.9
$ws.Range("I97").Value = 2593.1667
$ws.Range("J97").Value = 3475
$ws.Range("K97").Value = 2593.1667
$ws.Range("L97").Value = 3475
$ws.Range("M97").Value = -2097.1667
$ws.Range("N97").Value = -4467

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7343.8887
$ws.Range("J7").Value = 9459
$ws.Range("L7").Value = 9459
$ws.Range("N7").Value = -9683
$ws.Range("H126").Value = 7343.8887
$ws.Range("J126").Value = 9459
$ws.Range("L126").Value = 28377
$ws.Range("N126").Value = -33317
